# Update the "已预约/人数" (column F) counter values across the sheets,
# matching the regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 257
$ws1.Range("F3").Value  = 322
$ws1.Range("F17").Value = 7131
$ws1.Range("F28").Value = 1917
$ws1.Range("F33").Value = 291
$ws1.Range("F34").Value = 43
$ws1.Range("F36").Value = 1231
$ws1.Range("F37").Value = 2784
$ws1.Range("F41").Value = 397

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F12").Value = 299

# Sheet "本地生活" (index 3)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 78

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 322
$ws4.Range("F13").Value = 78
$ws4.Range("F19").Value = 7131
$ws4.Range("F25").Value = 299
$ws4.Range("F30").Value = 1917
$ws4.Range("F36").Value = 291
$ws4.Range("F37").Value = 43
$ws4.Range("F39").Value = 1231
$ws4.Range("F41").Value = 2784
$ws4.Range("F45").Value = 397
